# Generate Report for Handback
# Rebuild Overview / zh-cn / de-de sheets with latest handback status.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

function Set-CellText {
    param($ws, $row, $col, $text)
    $ws.Cells.Item($row, $col).Value2 = $text
}

function Add-Link {
    param($ws, $cellRef, $url, $displayText)
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, [System.Type]::Missing, [System.Type]::Missing, $displayText) | Out-Null
}

# ---------------- Overview sheet ----------------
$wsOverview.Hyperlinks.Delete()
Set-CellText $wsOverview 2 2 "Handed back: in sync with en-US"
Set-CellText $wsOverview 2 3 "Handed back: in sync with en-US"
Set-CellText $wsOverview 2 4 "2016-50-18 05:50:08"
Add-Link $wsOverview "A2" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsOverview 3 2 "Handed back: in sync with en-US"
Set-CellText $wsOverview 3 3 "Handed back: in sync with en-US"
Set-CellText $wsOverview 3 4 "2016-50-18 05:50:08"
Add-Link $wsOverview "A3" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md" "6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md"
Set-CellText $wsOverview 4 2 "Handed back: in sync with en-US"
Set-CellText $wsOverview 4 3 "Handed back: in sync with en-US"
Set-CellText $wsOverview 4 4 "2016-47-18 05:47:53"
Add-Link $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/047f45cf64a8750278f38d2243cf50fefe658bc7/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Set-CellText $wsOverview 5 2 "Handback transform failed"
Set-CellText $wsOverview 5 3 "Handback transform failed"
Set-CellText $wsOverview 5 4 "2016-48-18 05:48:46"
Add-Link $wsOverview "A5" "https://github.com/OpenLocalizationTest/oltest/blob/184fea56913b2dba6f94eed6dfca403f8b31e5f0/e2e/48b1637f-ef37-49a1-9da0-eb8cfcc87947.md" "48b1637f-ef37-49a1-9da0-eb8cfcc87947.md"

# ---------------- zh-cn sheet ----------------
$wsZhCn.Hyperlinks.Delete()
Set-CellText $wsZhCn 2 1 "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsZhCn 2 2 ".md"
Set-CellText $wsZhCn 2 3 "Handed back: in sync with en-US"
Set-CellText $wsZhCn 2 4 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"
Set-CellText $wsZhCn 2 5 "2016-03-18 05:50:05"
Set-CellText $wsZhCn 2 6 "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsZhCn 2 7 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"
Set-CellText $wsZhCn 2 8 "2016-03-18 05:50:21"
Set-CellText $wsZhCn 2 9 "Include"
$wsZhCn.Cells.Item(2, 11).ClearContents()
Add-Link $wsZhCn "A2" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Add-Link $wsZhCn "B2" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" ".md"
Add-Link $wsZhCn "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbb828def8f6bdac1e9749103cd781201223477b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf" "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"
Add-Link $wsZhCn "F2" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fbb828def8f6bdac1e9749103cd781201223477b/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Add-Link $wsZhCn "G2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fbb828def8f6bdac1e9749103cd781201223477b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf" "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"

Set-CellText $wsZhCn 3 1 "6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md"
Set-CellText $wsZhCn 3 2 ".md"
Set-CellText $wsZhCn 3 3 "Handed back: in sync with en-US"
Set-CellText $wsZhCn 3 4 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"
Set-CellText $wsZhCn 3 5 "2016-03-18 05:50:05"
Set-CellText $wsZhCn 3 6 "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsZhCn 3 7 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"
Set-CellText $wsZhCn 3 8 "2016-03-18 05:50:21"
Set-CellText $wsZhCn 3 9 "Include"
$wsZhCn.Cells.Item(3, 11).ClearContents()
Add-Link $wsZhCn "A3" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md" "6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md"
Add-Link $wsZhCn "B3" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md" ".md"
Add-Link $wsZhCn "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbb828def8f6bdac1e9749103cd781201223477b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf" "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"
Add-Link $wsZhCn "F3" $null "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Add-Link $wsZhCn "G3" $null "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.zh-cn.xlf"

Set-CellText $wsZhCn 4 1 "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Set-CellText $wsZhCn 4 2 ".md"
Set-CellText $wsZhCn 4 3 "Handed back: in sync with en-US"
Set-CellText $wsZhCn 4 4 "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf"
Set-CellText $wsZhCn 4 5 "2016-03-18 05:47:50"
Set-CellText $wsZhCn 4 6 "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Set-CellText $wsZhCn 4 7 "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf"
Set-CellText $wsZhCn 4 8 "2016-03-18 05:48:09"
Set-CellText $wsZhCn 4 9 "Include"
$wsZhCn.Cells.Item(4, 11).ClearContents()
Add-Link $wsZhCn "A4" "https://github.com/OpenLocalizationTest/oltest/blob/047f45cf64a8750278f38d2243cf50fefe658bc7/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Add-Link $wsZhCn "B4" "https://github.com/OpenLocalizationTest/oltest/blob/047f45cf64a8750278f38d2243cf50fefe658bc7/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" ".md"
Add-Link $wsZhCn "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/808d7184a05a85e37cb2c60c415c6e244e491a64/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf" "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf"
Add-Link $wsZhCn "F4" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fa9a877af2a8dc67a001995a162d254996f35493/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Add-Link $wsZhCn "G4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/57cfc8741a6e69966032a361d583510287489f6e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf" "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.zh-cn.xlf"

Set-CellText $wsZhCn 5 1 "48b1637f-ef37-49a1-9da0-eb8cfcc87947.md"
Set-CellText $wsZhCn 5 2 ".md"
Set-CellText $wsZhCn 5 3 "Handback transform failed"
Set-CellText $wsZhCn 5 4 "48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.zh-cn.xlf"
Set-CellText $wsZhCn 5 5 "2016-03-18 05:48:43"
$wsZhCn.Cells.Item(5, 6).ClearContents()
$wsZhCn.Cells.Item(5, 7).ClearContents()
Set-CellText $wsZhCn 5 8 "0001-01-01 00:00:00"
Set-CellText $wsZhCn 5 9 "Include"
Set-CellText $wsZhCn 5 11 "Handback file name: 1pfdtgmw.wwa is different with handoff file name: 48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.zh-cn."
Add-Link $wsZhCn "A5" "https://github.com/OpenLocalizationTest/oltest/blob/184fea56913b2dba6f94eed6dfca403f8b31e5f0/e2e/48b1637f-ef37-49a1-9da0-eb8cfcc87947.md" "48b1637f-ef37-49a1-9da0-eb8cfcc87947.md"
Add-Link $wsZhCn "B5" "https://github.com/OpenLocalizationTest/oltest/blob/184fea56913b2dba6f94eed6dfca403f8b31e5f0/e2e/48b1637f-ef37-49a1-9da0-eb8cfcc87947.md" ".md"
Add-Link $wsZhCn "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c119ca7f6b3ac4347e7190b014f6c0edec1d837/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.zh-cn.xlf" "48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.zh-cn.xlf"

# ---------------- de-de sheet ----------------
$wsDeDe.Hyperlinks.Delete()
Set-CellText $wsDeDe 2 1 "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsDeDe 2 2 ".md"
Set-CellText $wsDeDe 2 3 "Handed back: in sync with en-US"
Set-CellText $wsDeDe 2 4 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"
Set-CellText $wsDeDe 2 5 "2016-03-18 05:50:08"
Set-CellText $wsDeDe 2 6 "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsDeDe 2 7 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"
Set-CellText $wsDeDe 2 8 "2016-03-18 05:50:26"
Set-CellText $wsDeDe 2 9 "Include"
$wsDeDe.Cells.Item(2, 11).ClearContents()
Add-Link $wsDeDe "A2" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Add-Link $wsDeDe "B2" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" ".md"
Add-Link $wsDeDe "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20f6e7a6ab637a3ec47a14e0fd231b27ca48cdee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf" "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"
Add-Link $wsDeDe "F2" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/20f6e7a6ab637a3ec47a14e0fd231b27ca48cdee/e2e/02ae41b2-4940-4f0c-8272-815588b7e66c.md" "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Add-Link $wsDeDe "G2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/20f6e7a6ab637a3ec47a14e0fd231b27ca48cdee/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf" "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"

Set-CellText $wsDeDe 3 1 "6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md"
Set-CellText $wsDeDe 3 2 ".md"
Set-CellText $wsDeDe 3 3 "Handed back: in sync with en-US"
Set-CellText $wsDeDe 3 4 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"
Set-CellText $wsDeDe 3 5 "2016-03-18 05:50:08"
Set-CellText $wsDeDe 3 6 "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Set-CellText $wsDeDe 3 7 "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"
Set-CellText $wsDeDe 3 8 "2016-03-18 05:50:26"
Set-CellText $wsDeDe 3 9 "Include"
$wsDeDe.Cells.Item(3, 11).ClearContents()
Add-Link $wsDeDe "A3" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md" "6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md"
Add-Link $wsDeDe "B3" "https://github.com/OpenLocalizationTest/oltest/blob/7089534ac8b93c0410bd85e3a425d8ef7dd9bf3e/e2e/6448ebcc-5ff5-48a5-bbae-c682a603fe0e.md" ".md"
Add-Link $wsDeDe "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20f6e7a6ab637a3ec47a14e0fd231b27ca48cdee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf" "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"
Add-Link $wsDeDe "F3" $null "02ae41b2-4940-4f0c-8272-815588b7e66c.md"
Add-Link $wsDeDe "G3" $null "02ae41b2-4940-4f0c-8272-815588b7e66c.f66d89511b66a33246064de26263e70bedcf9f0c.de-de.xlf"

Set-CellText $wsDeDe 4 1 "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Set-CellText $wsDeDe 4 2 ".md"
Set-CellText $wsDeDe 4 3 "Handed back: in sync with en-US"
Set-CellText $wsDeDe 4 4 "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf"
Set-CellText $wsDeDe 4 5 "2016-03-18 05:47:53"
Set-CellText $wsDeDe 4 6 "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Set-CellText $wsDeDe 4 7 "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf"
Set-CellText $wsDeDe 4 8 "2016-03-18 05:48:13"
Set-CellText $wsDeDe 4 9 "Include"
$wsDeDe.Cells.Item(4, 11).ClearContents()
Add-Link $wsDeDe "A4" "https://github.com/OpenLocalizationTest/oltest/blob/047f45cf64a8750278f38d2243cf50fefe658bc7/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Add-Link $wsDeDe "B4" "https://github.com/OpenLocalizationTest/oltest/blob/047f45cf64a8750278f38d2243cf50fefe658bc7/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" ".md"
Add-Link $wsDeDe "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a6a50fb47ee94db8c8b9ce8ca8c9883a89b81f89/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf" "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf"
Add-Link $wsDeDe "F4" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/379978af6614c4f7d525a3caca7ab956ba298ae4/e2e/97814758-b34a-46ff-8abf-bdbd317417c5.md" "97814758-b34a-46ff-8abf-bdbd317417c5.md"
Add-Link $wsDeDe "G4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5b76307d17054f5abaae741ff280d76e4266d582/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf" "97814758-b34a-46ff-8abf-bdbd317417c5.16eb167fb079c913182a12952c834c1a7ba86fa0.de-de.xlf"

Set-CellText $wsDeDe 5 1 "48b1637f-ef37-49a1-9da0-eb8cfcc87947.md"
Set-CellText $wsDeDe 5 2 ".md"
Set-CellText $wsDeDe 5 3 "Handback transform failed"
Set-CellText $wsDeDe 5 4 "48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.de-de.xlf"
Set-CellText $wsDeDe 5 5 "2016-03-18 05:48:46"
$wsDeDe.Cells.Item(5, 6).ClearContents()
$wsDeDe.Cells.Item(5, 7).ClearContents()
Set-CellText $wsDeDe 5 8 "0001-01-01 00:00:00"
Set-CellText $wsDeDe 5 9 "Include"
Set-CellText $wsDeDe 5 11 "Handback file name: 1pfdtgmw.wwa is different with handoff file name: 48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.de-de."
Add-Link $wsDeDe "A5" "https://github.com/OpenLocalizationTest/oltest/blob/184fea56913b2dba6f94eed6dfca403f8b31e5f0/e2e/48b1637f-ef37-49a1-9da0-eb8cfcc87947.md" "48b1637f-ef37-49a1-9da0-eb8cfcc87947.md"
Add-Link $wsDeDe "B5" "https://github.com/OpenLocalizationTest/oltest/blob/184fea56913b2dba6f94eed6dfca403f8b31e5f0/e2e/48b1637f-ef37-49a1-9da0-eb8cfcc87947.md" ".md"
Add-Link $wsDeDe "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0dd8e902a8f29ca90719e4ec2698c6ed9af20098/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.de-de.xlf" "48b1637f-ef37-49a1-9da0-eb8cfcc87947.946552f7dc47191d65abcdd17529a318bada893d.de-de.xlf"

